# Update countries & provincias Spain
# Applies the data refresh described by the commit: updates the "last
# updated" timestamp, refreshes several per-country case counters, and
# re-labels four row-pairs whose countries swapped order in the source
# feed (Belgica/Mexico, Dinamarca/Corea del Sur, Santa Lucia/Belice,
# Montserrat/Groenlandia).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp banner (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 03:35"

# --- Row 16: China -------------------------------------------------------
$ws.Range("B16").Value = 82967
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 78249
$ws.Range("E16").Value = 84

# --- Rows 19-20: Belgica / Mexico swap places + refreshed counters ------
$ws.Range("A19").Value = "Mexico"
$ws.Range("B19").Value = 56594
$ws.Range("C19").Value = 2248
$ws.Range("D19").Value = 38876
$ws.Range("E19").Value = 11628
$ws.Range("G19").Value = 424
$ws.Range("H19").Value = 6090

$ws.Range("A20").Value = "Belgica"
$ws.Range("B20").Value = 55983
$ws.Range("D20").Value = 14847
$ws.Range("E20").Value = 31986
$ws.Range("H20").Value = 9150

# --- Rows 47-48: Dinamarca / Corea del Sur swap places -------------------
$ws.Range("A47").Value = "Corea del Sur"
$ws.Range("B47").Value = 11122
$ws.Range("C47").Value = 12
$ws.Range("D47").Value = 10135
$ws.Range("E47").Value = 723
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 264

$ws.Range("A48").Value = "Dinamarca"
$ws.Range("B48").Value = 11117
$ws.Range("D48").Value = 9536
$ws.Range("E48").Value = 1027
$ws.Range("H48").Value = 554

# --- Row 96: Nueva Zelanda ------------------------------------------------
$ws.Range("D96").Value = 1452
$ws.Range("E96").Value = 30

# --- Row 129: Jamaica ------------------------------------------------------
$ws.Range("B129").Value = 529
$ws.Range("C129").Value = 9
$ws.Range("D129").Value = 171
$ws.Range("E129").Value = 349

# --- Rows 197-198: Santa Lucia / Belice swap places -----------------------
$ws.Range("A197").Value = "Belice"
$ws.Range("D197").Value = 16
$ws.Range("H197").Value = 2

$ws.Range("A198").Value = "Santa Lucia"
$ws.Range("D198").Value = 18
$ws.Range("H198").Value = 0

# --- Rows 209-210: Montserrat / Groenlandia swap places -------------------
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 0

$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1
